$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Sep 11 13:43:10 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 13:43:23 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 13:43:35 EDT 2023"
$ws.Range("B5").Value = "Mon Sep 11 13:43:47 EDT 2023"
$ws.Range("B6").Value = "Mon Sep 11 13:44:00 EDT 2023"
$ws.Range("B7").Value = "Mon Sep 11 13:44:12 EDT 2023"
$ws.Range("B8").Value = "Mon Sep 11 13:44:24 EDT 2023"
$ws.Range("B9").Value = "Mon Sep 11 13:44:36 EDT 2023"
$ws.Range("B10").Value = "Mon Sep 11 13:44:48 EDT 2023"
$ws.Range("B11").Value = "Mon Sep 11 13:45:00 EDT 2023"
$ws.Range("B12").Value = "Mon Sep 11 13:45:12 EDT 2023"
$ws.Range("B13").Value = "Mon Sep 11 13:45:24 EDT 2023"
$ws.Range("B14").Value = "Mon Sep 11 13:45:36 EDT 2023"
$ws.Range("B15").Value = "Mon Sep 11 13:45:49 EDT 2023"
$ws.Range("B16").Value = "Mon Sep 11 13:46:01 EDT 2023"
$ws.Range("B17").Value = "Mon Sep 11 13:46:13 EDT 2023"
$ws.Range("B18").Value = "Mon Sep 11 13:46:24 EDT 2023"
$ws.Range("B19").Value = "Mon Sep 11 13:46:36 EDT 2023"
$ws.Range("B20").Value = "Mon Sep 11 13:46:48 EDT 2023"
$ws.Range("B21").Value = "Mon Sep 11 13:47:00 EDT 2023"
$ws.Range("B22").Value = "Mon Sep 11 13:47:12 EDT 2023"
$ws.Range("B23").Value = "Mon Sep 11 13:47:24 EDT 2023"
$ws.Range("B24").Value = "Mon Sep 11 13:47:35 EDT 2023"
$ws.Range("B25").Value = "Mon Sep 11 13:47:47 EDT 2023"
$ws.Range("B26").Value = "Mon Sep 11 13:47:59 EDT 2023"
$ws.Range("B27").Value = "Mon Sep 11 13:48:11 EDT 2023"
$ws.Range("B28").Value = "Mon Sep 11 13:48:23 EDT 2023"
$ws.Range("B29").Value = "Mon Sep 11 13:48:35 EDT 2023"
